$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a new data row at row 74 (pushes the existing rows 74.. down by
#    one, including the "Total" row and the footer row).
# ---------------------------------------------------------------------------
$ws.Rows.Item(74).Insert()

# Copy the formatting (styles / merges) of the row that is now directly
# below the freshly inserted blank row (row 75, which used to be row 74)
# so the new row 74 looks exactly like the other data rows.
$ws.Range("A75:Q75").Copy()
$ws.Range("A74:Q74").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# Re-create the merged cells for the new row (format paste does not restore
# merges).
$ws.Range("A74:B74").Merge()
$ws.Range("C74:G74").Merge()
$ws.Range("H74:K74").Merge()
$ws.Range("L74:M74").Merge()
$ws.Range("N74:O74").Merge()

# ---------------------------------------------------------------------------
# 2) Fill in the values for the new product row.
# ---------------------------------------------------------------------------
$ws.Range("A74").Value = 68
$ws.Range("C74").Value = "شامبو جونسون وسط"
$ws.Range("H74").Value = "1:0"

# L74 ("order limit") must stay text ("0") even though the cell's number
# format is numeric - force text storage the same way as for P74 below.
$fmtL = $ws.Range("L74").NumberFormat
$ws.Range("L74").NumberFormat = "@"
$ws.Range("L74").Value = "0"
$ws.Range("L74").NumberFormat = $fmtL

$ws.Range("N74").Value = "85.00"

# P74 ("sell price") must be stored as literal text "85.0000" although the
# cell's number format is numeric (0.00) - temporarily switch the cell to a
# text format, assign the value, then restore the original numeric format
# so the stored value keeps its text type while the visual style id is kept.
$fmtP = $ws.Range("P74").NumberFormat
$ws.Range("P74").NumberFormat = "@"
$ws.Range("P74").Value = "85.0000"
$ws.Range("P74").NumberFormat = $fmtP

$ws.Range("Q74").Value = "1:0"

# ---------------------------------------------------------------------------
# 3) Update the "Total" row (now row 83) with the new sum of the sell-price
#    column, and update the footer timestamp (now row 84).
# ---------------------------------------------------------------------------
$ws.Range("P83").Value = 3582.685

$ws.Range("A84").Value = "Saturday, 7 June, 2025 10:44 PM"

# ---------------------------------------------------------------------------
# 4) Fix up row heights: in this workbook row heights follow a fixed
#    pattern based on the absolute row number (25.5 / 24.75 / 25.5 / 24.75 /
#    25.5 repeating every 5 rows starting at row 7), and the footer row is
#    always 16.5. Native row-insert keeps heights attached to the shifted
#    content, so re-apply the expected heights explicitly.
# ---------------------------------------------------------------------------
$ws.Rows.Item(74).RowHeight = 25.5
$ws.Rows.Item(75).RowHeight = 24.75
$ws.Rows.Item(76).RowHeight = 25.5
$ws.Rows.Item(77).RowHeight = 25.5
$ws.Rows.Item(78).RowHeight = 24.75
$ws.Rows.Item(79).RowHeight = 25.5
$ws.Rows.Item(80).RowHeight = 24.75
$ws.Rows.Item(81).RowHeight = 25.5
$ws.Rows.Item(82).RowHeight = 25.5
$ws.Rows.Item(83).RowHeight = 24.75
$ws.Rows.Item(84).RowHeight = 16.5
